# ID 23 and associated to ticket #7:
# The height isn't changing anymore. Now a bubble appears at left of search
# field. This bubble shows the error when we put the mouse over the bubble.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Backlog")
$ws.Activate()

# --- Row 23: "to set a search background color..." feature replaced by
#     the new bubble feature, state moves to "To test", and Real(hours) now set
$ws.Range("D23").Value2 = "To test"
$ws.Range("E23").Value2 = "to set a bubble to advertise the users"
$ws.Range("F23").Value2 = "to see whether there is a error on search"
$ws.Range("H23").Value2 = 1

# --- Row 25: Real (hours) now filled in
$ws.Range("H25").Value2 = 1

# --- Row 28: previously held the "status bar" ticket, now replaced with what
#     used to be row 29's ticket content ("to save as a list of locked variable")
$ws.Range("E28").Value2 = "to save as a list of locked variable"
$ws.Range("F28").Value2 = "to load quicly a list of locked elements"
$ws.Range("G28").Value2 = 3

# --- Row 29: now empty (its former content moved up into row 28)
$ws.Range("B29:I29").ClearContents()

# --- Remove the now-unused trailing blank row 100
$ws.Rows.Item(100).Delete()

# --- Update conditional formatting ranges that used to extend to row 100
$fcs = $ws.Cells.FormatConditions
for ($i = 1; $i -le $fcs.Count; $i++) {
    $fc = $fcs.Item($i)
    $addr = $fc.AppliesTo.Address()
    if ($addr -like "*`$C`$3:`$C`$100*") {
        $fc.ModifyAppliesToRange($ws.Range("C3:C99"))
    } elseif ($addr -like "*`$D`$3:`$D`$100*") {
        $fc.ModifyAppliesToRange($ws.Range("D3:D99"))
    }
}

# --- Shrink the AutoFilter range from C2:D29 to C2:D28, keeping the filter
$ws.AutoFilterMode = $false
[void]$ws.Range("C2:D28").AutoFilter(2, @("To do", "To test"), 7)

# --- Update the hidden _FilterDatabase defined name to match the new filter range
for ($i = 1; $i -le $wb.Names.Count; $i++) {
    $n = $wb.Names.Item($i)
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "=Backlog!`$C`$2:`$D`$28"
    }
}

# --- Move the active selection to D23
[void]$ws.Range("D23").Select()
